# Automatic Word resume update
# Adds a new "Docker / Rancher container management" bullet immediately
# after the "Continuous integration" bullet in the "General skills" list.

$d = $word.ActiveDocument

# Locate the exact "Continuous integration" bullet (case-sensitive so we
# don't match the unrelated "Continuous Integration" heading later in the
# document, under the Bentley Systems experience entry).
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "Continuous integration",  # FindText
    $true,                     # MatchCase
    $true,                     # MatchWholeWord
    $false,                    # MatchWildcards
    $false,                    # MatchSoundsLike
    $false,                    # MatchAllWordForms
    $true,                     # Forward
    1,                         # Wrap (wdFindContinue)
    $false,                    # Format
    "",                        # ReplaceWith
    0                          # Replace (wdReplaceNone)
)

if (-not $found) {
    throw "Could not find the 'Continuous integration' bullet to anchor the new entry."
}

# The found range corresponds to exactly one paragraph; grab that paragraph.
$anchorPara = $searchRange.Paragraphs(1)

# Build the new list-item paragraph as raw OOXML so it carries the same
# paragraph style ("Compact") and numbering (ilvl 0 / numId 1002) as its
# sibling bullets, and keep the xml:space="preserve" markup on the text run
# consistent with the rest of the document.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1002"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Docker / Rancher container management</w:t></w:r></w:p>'

# Inserting this XML "at" the anchor paragraph's range appends the new
# paragraph immediately after it (without disturbing the anchor paragraph's
# own content). Cast to [void] so the method's return value isn't echoed.
[void]$anchorPara.Range.InsertXML($newParaXml, 0)
